$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 124
$ws1.Range("F5").Value = 2821
$ws1.Range("F6").Value = 772
$ws1.Range("F9").Value = 76
$ws1.Range("F12").Value = 424
$ws1.Range("F15").Value = 1247
$ws1.Range("F18").Value = 2656
$ws1.Range("F24").Value = 559
$ws1.Range("F32").Value = 4637

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 18
$ws2.Range("F8").Value = 344
$ws2.Range("F12").Value = 163
$ws2.Range("F21").Value = 271
$ws2.Range("F32").Value = 497

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 197

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 124
$ws4.Range("F6").Value = 197
$ws4.Range("F10").Value = 2821
$ws4.Range("F11").Value = 772
$ws4.Range("F14").Value = 76
$ws4.Range("F17").Value = 424
$ws4.Range("F19").Value = 344
$ws4.Range("F23").Value = 1247
$ws4.Range("F27").Value = 2656
$ws4.Range("F37").Value = 559
$ws4.Range("F38").Value = 559
$ws4.Range("F39").Value = 271
$ws4.Range("F47").Value = 4637
$ws4.Range("F50").Value = 497
